$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.116.97"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "3.148.47"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.52"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.140.25"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("E12").Value = "  -1.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.26"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("D15").Value = "3.670.27"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.28"
$ws.Range("E17").Value = "  +1.61%  "

$ws.Range("D18").Value = "63.943.95"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "3.148.26"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.66"
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.44"
$ws.Range("E25").Value = "  -1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +6.04%  "

$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("E28").Value = "  +7.88%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("E31").Value = "  +7.37%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.72"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +2.08%  "

$ws.Range("D35").Value = "0.0₃0838"
$ws.Range("E35").Value = "  -4.26%  "

$ws.Range("E36").Value = "  +1.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.18"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("E38").Value = "  -2.73%  "

$ws.Range("E39").Value = "  -5.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "461.77"
$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.32"
$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.23"
$ws.Range("E42").Value = "  +5.60%  "

$ws.Range("E43").Value = "  +4.97%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").Value = "2.930.27"
$ws.Range("E45").Value = "  +0.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.66"
$ws.Range("E46").Value = "  +11.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.109"
$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.31"
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("E50").Value = "  +2.79%  "

